$d = $word.ActiveDocument

# Curly right-single-quote used throughout the document for "Simpson's".
$rsquo = [char]0x2019
# Temporary marker character (removed below) used only to locate the exact
# character offset where the "_GoBack" bookmark needs to live once the
# replacement text is in place.
$marker = [char]0x0007

# Original sentence (spans several runs, but Find/Replace matches across
# run boundaries just fine):
#   "; it does poorly when it comes to approximating very curvy equations.
#    Simpson's Rule, however, "
# becomes two sentences split into two paragraphs:
#   "; as a result, the Composite Trapezoid Rule performs poorly when it
#    comes to approximating very curvy equations, like ours which involves
#    a few sine and cosines. "
#   (new paragraph) "Simpson's Rule, however, "...
$old = "; it does poorly when it comes to approximating very curvy equations. Simpson" + $rsquo + "s Rule, however, "
$new = "; as a result, the Composite Trapezoid Rule performs poorly when it comes to approxima" + $marker + "ting very curvy equations, like ours which involves a few sine and cosines. ^pSimpson" + $rsquo + "s Rule, however, "

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Re-home the document's "_GoBack" bookmark onto the marker's position (this
# mirrors Word automatically relocating _GoBack to the site of the most
# recent edit), then remove the one-character marker itself.
$markerRange = $d.Content
$null = $markerRange.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bmRange = $d.Range($markerRange.Start, $markerRange.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
$markerRange.Text = ""
